$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header row (row 1) text:
#    columns A..J  (1..10): "<name>_old" -> "<name>_FV2304"
#    column  K     (11)   : "diff" stays unchanged
#    columns L..U  (12..21): "<name>_new" -> "<name>_FV2310"
for ($c = 1; $c -le 10; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $text = [string]$cell.Value2
    if ($text.EndsWith("_old")) {
        $cell.Value2 = $text.Substring(0, $text.Length - 4) + "_FV2304"
    }
}
for ($c = 12; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $text = [string]$cell.Value2
    if ($text.EndsWith("_new")) {
        $cell.Value2 = $text.Substring(0, $text.Length - 4) + "_FV2310"
    }
}

# 2. Turn the used range A1:U73 into an Excel table ("Table1") with a
#    header row, autofilter and banded rows.
$rng = $ws.Range("A1:U73")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# 3. Freeze the header row (View > Freeze Panes > Freeze Top Row):
#    select A2 so the split sits right below row 1, then freeze.
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

Write-Host "done"
